# Add a "2022-Q3" sheet of fund holdings, positioned between the existing
# "总计" and "2020-Q4" sheets, and push the previous "总计" summary row for
# 2020-Q4 down to make room for a new 2022-Q3 summary row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: shift the old 2020-Q4 row down to row 3, and put the
#    new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets("总计")

# Duplicate row 2's formatting into row 3 before we overwrite row 2.
$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2020-Q4"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.16

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.37

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet, inserted right before "2020-Q4".
# ---------------------------------------------------------------------
$before = $wb.Worksheets("2020-Q4")
$ws = $wb.Worksheets.Add($before)
$ws.Name = "2022-Q3"

# Header row: reuse the bold/bordered header style already used on "总计".
$summary.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Column A style (row index), reused from the summary sheet.
$summary.Range("A2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("H2").Value = 9
$ws.Range("H3").Value = 9

# B2:G3 all hold numeric-looking text (fund codes with leading zeros, and
# decimal figures that must keep trailing zeros) - force text storage for
# the whole block in one shot so only a single throwaway style is created.
$data = $ws.Range("B2:G3")
$data.NumberFormat = "@"
$ws.Range("B2").Value = "002446"
$ws.Range("C2").Value = "广发利鑫灵活配置混合A"
$ws.Range("D2").Value = "13.77"
$ws.Range("E2").Value = "74.30"
$ws.Range("F2").Value = "2.45"
$ws.Range("G2").Value = "0.3374"
$ws.Range("B3").Value = "011172"
$ws.Range("C3").Value = "广发利鑫灵活配置混合C"
$ws.Range("D3").Value = "1.41"
$ws.Range("E3").Value = "74.30"
$ws.Range("F3").Value = "2.45"
$ws.Range("G3").Value = "0.0345"
$data.ClearFormats()
